$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Header timestamp
$ws.Range("A1").Value = "Datos actualizados a 23 de Junio de 2020 a las 09:34"

# --- India (row 7) ---
$ws.Range("B7").Value = 441070
$ws.Range("C7").Value = 620
$ws.Range("D7").Value = 248286
$ws.Range("E7").Value = 178765
$ws.Range("G7").Value = 4
$ws.Range("H7").Value = 14019

# --- Singapur (row 35) ---
$ws.Range("B35").Value = 42432
$ws.Range("C35").Value = 119
$ws.Range("E35").Value = 6816

# --- Ucrania (row 38) ---
$ws.Range("B38").Value = 38074
$ws.Range("C38").Value = 833
$ws.Range("D38").Value = 16956
$ws.Range("E38").Value = 20083
$ws.Range("G38").Value = 23
$ws.Range("H38").Value = 1035

# --- Armenia moves above Nigeria (rows 52/53 swap with updated Armenia stats) ---
$ws.Range("A52").Value = "Armenia"
$ws.Range("B52").Value = 21006
$ws.Range("C52").Value = 418
$ws.Range("D52").Value = 10144
$ws.Range("E52").Value = 10490
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 12
$ws.Range("H52").Value = 372

$ws.Range("A53").Value = "Nigeria"
$ws.Range("B53").Value = 20919
$ws.Range("C53").Value = 0
$ws.Range("D53").Value = 7109
$ws.Range("E53").Value = 13285
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 525

# --- Chequia (row 67) ---
$ws.Range("B67").Value = 10561
$ws.Range("C67").Value = 38
$ws.Range("D67").Value = 7543
$ws.Range("E67").Value = 2681
$ws.Range("G67").Value = 1
$ws.Range("H67").Value = 337

# --- Estonia (row 106) ---
$ws.Range("B106").Value = 1982
$ws.Range("C106").Value = 1
$ws.Range("D106").Value = 1771
$ws.Range("E106").Value = 142

# --- Sri Lanka (row 108) ---
$ws.Range("D108").Value = 1548
$ws.Range("E108").Value = 392

# --- Lituania (row 112) ---
$ws.Range("B112").Value = 1803
$ws.Range("C112").Value = 2
$ws.Range("D112").Value = 1483
$ws.Range("E112").Value = 243
$ws.Range("G112").Value = 1
$ws.Range("H112").Value = 77

# --- Georgia (row 133) ---
$ws.Range("B133").Value = 911
$ws.Range("C133").Value = 3
$ws.Range("E133").Value = 136

# --- Montenegro (row 156) ---
$ws.Range("B156").Value = 375
$ws.Range("C156").Value = 8
$ws.Range("E156").Value = 51

# --- Siria (row 164) ---
$ws.Range("D164").Value = 86
$ws.Range("E164").Value = 126
